$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 177: revised figures + new trailing columns P:T ---
$ws.Range("D177").Value = 7187.2
$ws.Range("F177").Value = 6743.5
$ws.Range("G177").Value = 15932.1
$ws.Range("K177").Value = 72170.39999999999
$ws.Range("L177").Value = 4028.3
$ws.Range("O177").Value = 1474.7
$ws.Range("P177").Value = 76003
$ws.Range("Q177").Value = 117868.6
$ws.Range("R177").Value = 1327.1
$ws.Range("S177").Value = 52283.1
$ws.Range("T177").Value = 17712.5

# --- New rows 178-182: column A holds the date label (forced to Text so it
#     is not auto-converted to a date serial number) ---
$dateRng = $ws.Range("A178:A182")
$dateRng.NumberFormat = "@"
$ws.Range("A178").Value = "07-09-2021"
$ws.Range("A179").Value = "08-09-2021"
$ws.Range("A180").Value = "09-09-2021"
$ws.Range("A181").Value = "10-09-2021"
$ws.Range("A182").Value = "13-09-2021"
$dateRng.Style = "Normal"

# --- Row 178 data ---
$ws.Range("B178").Value = 35100
$ws.Range("C178").Value = 15374.3
$ws.Range("D178").Value = 7149.4
$ws.Range("E178").Value = 29916.1
$ws.Range("F178").Value = 6726.1
$ws.Range("G178").Value = 15843.1
$ws.Range("H178").Value = 3187.4
$ws.Range("I178").Value = 4992.8
$ws.Range("J178").Value = 1583.5
$ws.Range("K178").Value = 71799
$ws.Range("L178").Value = 4010.2
$ws.Range("M178").Value = 1636.5
$ws.Range("N178").Value = 17428.9
$ws.Range("O178").Value = 1454.7
$ws.Range("P178").Value = 79004.8
$ws.Range("R178").Value = 1321
$ws.Range("S178").Value = 51950.5
$ws.Range("T178").Value = 17598.4

# --- Row 179 data ---
$ws.Range("B179").Value = 35031.1
$ws.Range("C179").Value = 15286.6
$ws.Range("D179").Value = 7095.5
$ws.Range("E179").Value = 30181.2
$ws.Range("F179").Value = 6668.9
$ws.Range("G179").Value = 15610.3
$ws.Range("H179").Value = 3163
$ws.Range("I179").Value = 4972.1
$ws.Range("J179").Value = 1597.6
$ws.Range("K179").Value = 70970.39999999999
$ws.Range("L179").Value = 4017.4
$ws.Range("M179").Value = 1640.5
$ws.Range("N179").Value = 17270.5
$ws.Range("O179").Value = 1433.9
$ws.Range("P179").Value = 77459
$ws.Range("Q179").Value = 113412.8
$ws.Range("R179").Value = 1320.4
$ws.Range("S179").Value = 51469.7
$ws.Range("T179").Value = 17410.2

# --- Row 180 data ---
$ws.Range("B180").Value = 34879.4
$ws.Range("C180").Value = 15248.3
$ws.Range("D180").Value = 7024.2
$ws.Range("E180").Value = 30008.2
$ws.Range("F180").Value = 6684.7
$ws.Range("G180").Value = 15623.2
$ws.Range("H180").Value = 3114.7
$ws.Range("I180").Value = 4970
$ws.Range("J180").Value = 1578.9
$ws.Range("K180").Value = 70980.2
$ws.Range("L180").Value = 3993.6
$ws.Range("M180").Value = 1629.1
$ws.Range("N180").Value = 17304.3
$ws.Range("O180").Value = 1452.7
$ws.Range("P180").Value = 77127
$ws.Range("Q180").Value = 115360.9
$ws.Range("R180").Value = 1322.8
$ws.Range("S180").Value = 51395.1
$ws.Range("T180").Value = 17550.8

# --- Row 181 data ---
$ws.Range("B181").Value = 34607.7
$ws.Range("C181").Value = 15115.5
$ws.Range("D181").Value = 7029.2
$ws.Range("E181").Value = 30381.8
$ws.Range("F181").Value = 6663.8
$ws.Range("G181").Value = 15609.8
$ws.Range("H181").Value = 3125.8
$ws.Range("I181").Value = 5013.5
$ws.Range("J181").Value = 1576
$ws.Range("K181").Value = 71091.2
$ws.Range("L181").Value = 4002.7
$ws.Range("M181").Value = 1635.4
$ws.Range("N181").Value = 17474.6
$ws.Range("O181").Value = 1438.7
$ws.Range("P181").Value = 75895
$ws.Range("Q181").Value = 114285.9
$ws.Range("R181").Value = 1320.5
$ws.Range("S181").Value = 51521.8
$ws.Range("T181").Value = 17715.5

# --- Row 182 data ---
$ws.Range("D182").Value = 7083.1
$ws.Range("E182").Value = 30447.4
$ws.Range("F182").Value = 6721
$ws.Range("G182").Value = 15771.3
$ws.Range("H182").Value = 3127.9
$ws.Range("I182").Value = 4991.7
$ws.Range("J182").Value = 1570.1
$ws.Range("K182").Value = 71569.3
$ws.Range("L182").Value = 4021.3
$ws.Range("M182").Value = 1633.8
$ws.Range("N182").Value = 17446.3
$ws.Range("O182").Value = 1443.8
